$d = $word.ActiveDocument

# --- Locate the target paragraph: the one ending in
# "...(Responsabilidad de todo el equipo)." inside the retrospective table.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Responsabilidad de todo el equipo*") {
        $targetIndex = $i
        break
    }
}

$para = $d.Paragraphs.Item($targetIndex)
$paraRange = $para.Range
$beforeEnd = $paraRange.End

# Append the new sentence (with a leading space) right after the existing
# text, still inside the same paragraph. A trailing sentinel character is
# appended too so that, further down, we can carve out an unambiguous
# in-text position for the _GoBack bookmark (raw offsets that sit exactly
# on a paragraph boundary get mis-resolved to the next paragraph).
$newText = " Referente: Paula Pedrosa"
$sentinel = "\u0001"
$paraRange.InsertAfter($newText + $sentinel)

# Force the appended text into its own run (distinct w:r) even though its
# formatting ends up identical to the preceding run: toggle a character
# property on just the inserted span so the engine splits the run, then
# reset it back to match the original formatting.
$newSpan = $d.Range($beforeEnd - 1, $beforeEnd + $newText.Length)
$newSpan.Bold = 1
$newSpan.Bold = 0

# Remove the _GoBack bookmark from its old location later in the document.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# Re-create _GoBack collapsed right after "Paula Pedrosa" (i.e. right before
# the sentinel char, which is still safely inside the paragraph's text).
$sentinelPos = $beforeEnd - 1 + $newText.Length
$bkRange = $d.Range($sentinelPos, $sentinelPos)
$d.Bookmarks.Add("_GoBack", $bkRange)

# Drop the sentinel character now that the bookmark is anchored.
$d.Range($sentinelPos, $sentinelPos + 1).Text = ""
